$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.371.56'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '3.587.78'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.70'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.60'
$ws.Range('E6').Value = '  +2.98%  '
$ws.Range('D7').Value = '3.587.54'
$ws.Range('E7').Value = '  +1.15%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.87'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '4.193.62'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000205'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.52'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '3.593.93'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('E17').Value = '  +1.60%  '
$ws.Range('D18').Value = '66.392.97'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.12'
$ws.Range('E19').Value = '  -1.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.33'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.83'
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '422.01'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.611'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.25'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000120'
$ws.Range('E26').Value = '  +3.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.21'
$ws.Range('E27').Value = '  +5.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.36'
$ws.Range('E28').Value = '  +3.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.49'
$ws.Range('E29').Value = '  +0.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').Value = '3.583.65'
$ws.Range('E31').Value = '  +1.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.157'
$ws.Range('E32').Value = '  +4.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.00'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.41'
$ws.Range('E34').Value = '  -1.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.73'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.55'
$ws.Range('E37').Value = '  +0.71%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.66'
$ws.Range('E38').Value = '  -2.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '174.57'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0850'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.18'
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.879'
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '46.02'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('E44').Value = '  -2.58%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.50'
$ws.Range('E46').Value = '  +5.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.64'
$ws.Range('E47').Value = '  +4.01%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.21'
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.12'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.945'
$ws.Range('E51').Value = '  +2.05%  '
